$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column C ("Förändrad") from 45184 to 45186 for rows 2-12
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 3).Value2 = 45186
}

# 2) Add a friendly display-text second argument to the HYPERLINK formulas
#    in columns S, T, V, W, X, Y for rows 2-5 (the label equals column A's text).
$hyperlinkCols = @("S", "T", "V", "W", "X", "Y")
$hyperlinkPaths = @{
    "S" = "artfynd"
    "T" = "kartor"
    "V" = "klagomål"
    "W" = "klagomålsmail"
    "X" = "tillsyn"
    "Y" = "tillsynsmail"
}
$hyperlinkExt = @{
    "S" = "xlsx"
    "T" = "png"
    "V" = "docx"
    "W" = "docx"
    "X" = "docx"
    "Y" = "docx"
}

for ($r = 2; $r -le 5; $r++) {
    $label = $ws.Range("A$r").Value2
    foreach ($col in $hyperlinkCols) {
        $path = $hyperlinkPaths[$col]
        $ext = $hyperlinkExt[$col]
        $url = "https://klasma.github.io/Logging_HELSINGBORG/$path/$label.$ext"
        $formula = '=HYPERLINK("' + $url + '", "' + $label + '")'
        $ws.Range("$col$r").Formula = $formula
    }
}
